# Agent/entities.xlsx refactor:
#  - trim trailing whitespace from brand/model display names
#  - convert synonyms cells from Python-list-literal text to newline-joined text
#  - center+middle-align the whole table, wrap the synonyms column
#  - widen column C, grow rows 2/12-15/20 to fit wrapped text
#  - add a new "quantity" entity block (rows 18-20)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- trim trailing spaces on display-name cells ----
$ws.Range("A7").Value = "smartphoneName"
$ws.Range("B4").Value = "LG"
$ws.Range("B5").Value = "Sony  "
$ws.Range("B6").Value = "Google"
$ws.Range("B7").Value = "S9"
$ws.Range("B8").Value = "S8"
$ws.Range("B9").Value = "Pixel 2"

# ---- synonyms column: drop the ['a', 'b'] python-list formatting ----
$ws.Range("C2").Value = "smartphone`nphone`nmóvil`nterminal`nmóvil`nteléfono`nsmartphones"
$ws.Range("C3").Value = "samsung"
$ws.Range("C4").Value = "lg"
$ws.Range("C5").Value = "sony"
$ws.Range("C6").Value = "google"
$ws.Range("C7").Value = "s9"
$ws.Range("C8").Value = "s8"
$ws.Range("C9").Value = "pixel 2"
$ws.Range("C10").Value = "xperia xz1 compact"
$ws.Range("C11").Value = "pixel 2"
$ws.Range("C12").Value = "gama baja`nbaja"
$ws.Range("C13").Value = "gama media`nmedia"
$ws.Range("C14").Value = "gama alta`nalta"
$ws.Range("C15").Value = "gama premium`npremium"
$ws.Range("C16").Value = "test testing"
$ws.Range("C17").Value = "value"

# ---- new "quantity" entity rows ----
$ws.Range("A18").Value = "quantity"
$ws.Range("B18").Value = 3
$ws.Range("C18").Value = 3
$ws.Range("B19").Value = 4
$ws.Range("C19").Value = 4
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "el`ncuál`n1"
$ws.Range("A18:A20").Merge()

# ---- whole-table alignment: center horizontally + vertically ----
$ws.Range("A1:A20").HorizontalAlignment = -4108
$ws.Range("A1:A20").VerticalAlignment = -4108
$ws.Range("B1:C20").HorizontalAlignment = -4108
$ws.Range("B1:C20").VerticalAlignment = -4108

# ---- wrap + taller rows for the multi-line synonyms cells ----
$ws.Range("C2").WrapText = $true
$ws.Range("C12:C15").WrapText = $true
$ws.Range("C20").WrapText = $true

$ws.Rows.Item(2).RowHeight = 105
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 45

# ---- column C widened (69.15 "characters" settles to a stored width of 70) ----
$ws.Columns.Item(3).ColumnWidth = 69.15

# ---- selection parity with the authored workbook ----
$ws.Range("C26").Select()
